$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.264.69"
$ws.Range("E2").Value = "  -4.88%  "
$ws.Range("D3").Value = "2.609.22"
$ws.Range("E3").Value = "  -4.38%  "
$ws.Range("E4").Value = "  +0.16%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "519.24"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.99%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "141.56"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.12%  "
$ws.Range("E7").Value = "  +0.30%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "6.52"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -9.26%  "
$ws.Range("E10").Value = "  -3.94%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.334"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.31%  "
$ws.Range("E12").Value = "  +0.74%  "
$ws.Range("D13").Value = "3.068.05"
$ws.Range("E13").Value = "  -4.31%  "
$ws.Range("D14").Value = "58.243.60"
$ws.Range("E14").Value = "  -4.82%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "20.80"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.70%  "
$ws.Range("D16").Value = "2.657.31"
$ws.Range("E16").Value = "  -4.56%  "
$ws.Range("E17").Value = "  -2.87%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "335.90"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -3.36%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.38"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.47%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.33"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.13%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.25"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.67%  "
$ws.Range("E22").Value = "  +0.29%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "64.79"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.66%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.412"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.36%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.165"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -4.04%  "
$ws.Range("E26").Value = "  +0.89%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.10"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.68%  "
$ws.Range("D28").Value = "0.0₃0782"
$ws.Range("E28").Value = "  -5.89%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.53"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.11%  "
$ws.Range("E30").Value = "  +0.18%  "
$ws.Range("E31").Value = "  -1.44%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "18.70"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.45%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "149.92"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.19%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.06"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -5.42%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.17"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -5.39%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.888"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.04%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.844"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -7.04%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "36.14"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.01%  "
$ws.Range("E39").Value = "  -7.81%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.60"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.89%  "
$ws.Range("E41").Value = "  +0.15%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.602"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.84%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0965"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.82%  "
$ws.Range("B44").Value = "WhiteBITCoin"
$ws.Range("C44").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "10.64"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.04%  "
$ws.Range("B45").Value = "Bittensor"
$ws.Range("C45").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "266.54"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -6.58%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "19.02"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -7.46%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0526"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.79%  "
$ws.Range("D48").Value = "2.020.26"
$ws.Range("E48").Value = "  -5.40%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0227"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.94%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "4.55"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -8.65%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "18.11"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -7.07%  "
